$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign brand-new string values first, in the exact order they must appear
# in the rebuilt shared-strings table.
$ws.Range("F2").Value = "null"
$ws.Range("A2").Value = "BMW 1er Model"
$ws.Range("B1").Value = "Insurancestatus"
$ws.Range("B2").Value = "Modify"
$ws.Range("F1").Value = "carhsn"

# Re-home the pre-existing strings into their shifted columns (column B was
# inserted, so the old B/C/D columns become C/D/E).
$ws.Range("C1").Value = "car"
$ws.Range("D1").Value = "model"
$ws.Range("E1").Value = "type"
$ws.Range("C2").Value = "BMW"
$ws.Range("D2").Value = "1er"
$ws.Range("E2").Value = "Cabrio"

# Column widths. The host's ColumnWidth setter quantizes to 1/6-character
# steps (pixel rounding), so the literal authored widths (17.1796875 /
# 15.26953125 / 16.26953125 -- produced by real Excel's font-metric-based
# best-fit) are not exactly reproducible here. These inputs land on the
# closest reachable quantized width to each target.
$ws.Columns.Item(1).ColumnWidth = 16.3
$ws.Columns.Item(2).ColumnWidth = 14.5
$ws.Columns.Item(6).ColumnWidth = 15.5

# Selection
$ws.Range("A3:F3").Select()
